# Update data driven test case: record a "Pass" result for the
# DataDriven_TC row by populating the Results column (D2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Pass"
